$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace row 3 (UTICA UNIVERSITY) with LEE UNIVERSITY's data (previously in row 4).
# Numeric-looking values are prefixed with a leading apostrophe so Excel stores them
# as literal text (preserving thousands separators etc.), matching the source data,
# instead of silently converting them to numbers.
$ws.Range("A3").Value = "LEE UNIVERSITY"
$ws.Range("B3").Value = "https://projects.propublica.org/nonprofits/organizations/620502739/202401319349302970/full"
$ws.Range("C3").Value = "'97,254,831"
$ws.Range("D3").Value = "'108,753,431"
$ws.Range("E3").Value = "'93,462,606"
$ws.Range("F3").Value = "'105,914,787"
$ws.Range("G3").Value = "'166,167,918"
$ws.Range("H3").Value = "'160,791,459"
$ws.Range("I3").Value = "'0"
$ws.Range("J3").Value = "'10,985,421"
$ws.Range("K3").Value = "'0"
$ws.Range("L3").Value = "'30,946,403"
$ws.Range("M3").Value = "'0"
$ws.Range("N3").Value = "'41,931,824"
$ws.Range("O3").Value = "'0"
$ws.Range("P3").Value = "'6,600,350"
$ws.Range("Q3").Value = "'0"
$ws.Range("R3").Value = "'32,018,540"
$ws.Range("S3").Value = "'0"
$ws.Range("T3").Value = "'38,618,890"
$ws.Range("U3").Value = "'3.9"
$ws.Range("V3").Value = "'2.61"
$ws.Range("W3").Value = "'44.86"
$ws.Range("X3").Value = "'36.46"

# Remove the now-duplicate row 4 (old LEE UNIVERSITY row), shrinking the used range to A1:X3
$ws.Rows.Item(4).Delete()
